$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(1.189095389782267,1.091575191898755,1.031544212749566,1.007044431903353,1.002974102015429,1.03121394614692,1.155503292974117,1.397949951198996,1.575218013355993,1.65566124323567,1.686093120592716,1.679540445959276,1.658165508803961,1.645068749583231,1.569956732881735,1.523826144251814,1.49727462963466,1.488281637981459,1.528738744504551,1.664444681511782,1.752959290818751,1.705734206581553,1.526517852551024,1.332506991025241)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 2).Value = $values[$i]
}

$values = @(0.3143629987666365,0.3167364189788628,0.3183050860984693,0.3189724177096274,0.3190849263647095,0.3183139721355275,0.3151582857771302,0.3098501755620049,0.3064820801180019,0.3050643024824069,0.3045437954121653,0.3046551688331363,0.3050211523036594,0.3052474580102142,0.3065770306200903,0.3074219245033945,0.3079186560500986,0.308088692998254,0.3073308698938959,0.3049132102741865,0.3034285393823808,0.3042122305271491,0.3073720014419479,0.3111924296084325)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 4).Value = $values[$i]
}

$values = @(0.3146415609947795,0.3117002830617857,0.3099865883021877,0.3093116537425971,0.3092010014796926,0.3099773908330548,0.3136083834267254,0.321452205533447,0.3276454593973384,0.3305542151356562,0.3316686352040179,0.331428052795907,0.3306456410437022,0.3301680699524354,0.3274571870385756,0.3258174129267246,0.3248828747929053,0.3245679425040819,0.3259910794090786,0.3308751048700032,0.3341424118675889,0.3323917655340125,0.3259125392657793,0.3192539003107342)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 5).Value = $values[$i]
}

$values = @(0.9331639326599941,0.9435632741459656,0.9506257471101165,0.9536739018624303,0.9541903170358665,0.9506661669111907,0.9366089716705659,0.9144236209834489,0.9014144915596418,0.8962129804518568,0.8943464986254241,0.8947438851092144,0.8960573541032844,0.8968753396275915,0.9017688555558152,0.9049545004989241,0.9068542190848348,0.907509004482705,0.9046084040026869,0.8956687535512415,0.8904278967860648,0.8931699179284607,0.9047646616530756,0.9198482371428938)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 6).Value = $values[$i]
}

$values = @(0.4189333636495363,0.418791464117831,0.4190927542983616,0.4193130000293337,0.4193554524876149,0.4190953302374893,0.4188036974664016,0.42132358427952,0.4250744369049784,0.427196435899603,0.428059993722556,0.4278713389235946,0.4272662775050975,0.426903480640604,0.4249441437134038,0.4238487887640474,0.4232578856642419,0.4230645274134872,0.4239613411743761,0.4274423684837814,0.4300672654496083,0.4286342246488459,0.4239103352909552,0.4203093679110737)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 7).Value = $values[$i]
}

$values = @(0.5713768830932082,0.575501370408837,0.5783588284838146,0.5796049879978185,0.5798168475531753,0.578375303742277,0.5727315489807694,0.5642433541287488,0.5595803767529759,0.5578009269251112,0.5571762548960209,0.5573086021334035,0.5577485492772922,0.5580244332460751,0.5597035470050571,0.5608211782045629,0.5614961770771743,0.561730244126025,0.5606988751331698,0.5576179915852038,0.5558910488541073,0.5567865202869342,0.5607540672483253,0.5662634151716901)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 8).Value = $values[$i]
}

$values = @(0.8729593133135545,0.8848543325973672,0.8927819860345423,0.8961695223332313,0.8967415006239925,0.892827035997275,0.8769312070227855,0.8507102664761028,0.8344629819371612,0.8277267696119424,0.8252701039036765,0.8257950009696771,0.8275227702712087,0.8285933465956035,0.8349163878528074,0.8389630942569646,0.8413522792258732,0.8421717995127196,0.8385259370782236,0.8270127258866253,0.820037218569901,0.8237099271436819,0.8387233805994327,0.8572738205066486)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 9).Value = $values[$i]
}

$values = @(0.4129581985981048,0.4021598632990475,0.3956205775573665,0.392979055650386,0.3925418521368869,0.3955848582041597,0.4092163079166653,0.4366521172464957,0.457218368457518,0.4666592210741811,0.4702460867903682,0.4694730721418239,0.4669540802184429,0.4654126516228416,0.4566030636605944,0.4512201740205484,0.4481321193259049,0.4470879511881662,0.4517923626948601,0.467693651801369,0.4781547628931975,0.472565324845732,0.451533655478201,0.4291566956025719)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 10).Value = $values[$i]
}

$values = @(0.4462141630035603,0.4003870681464434,0.3721440983002822,0.360609405164297,0.3586925685245035,0.371988639208638,0.4304353083547596,0.5441808231879861,0.627180920667854,0.6648085767201337,0.6790377528964768,0.6759741324324295,0.6659796180026945,0.659855107884141,0.6247191750920535,0.6031305505356954,0.5907012130027454,0.5864907991923758,0.6054299591706922,0.6689157888644388,0.710292753913734,0.6882199043878643,0.6043904517423755,0.5135067084843854)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 12).Value = $values[$i]
}

$values = @(0.3644955490078203,0.3312478445537153,0.3107695632379333,0.3024090415616527,0.3010198691089201,0.310656872030755,0.3530454144830415,0.4356364719167232,0.4959645293512835,0.5233277551934776,0.5336773886191111,0.5314489643749454,0.5241794746576858,0.5197250900676096,0.4941746040506843,0.4784791547219811,0.4694440283288088,0.4663836191329409,0.4801507452234546,0.5263150371916225,0.5564144962512501,0.5403566198139487,0.4793950552497677,0.4133532348656033)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 13).Value = $values[$i]
}

$values = @(1.929592594718571,1.93763772178923,1.944071012645551,1.947067751856167,1.947588001418282,1.944109909461702,1.932056297837093,1.920291942675959,1.918920134439873,1.919882319269334,1.920475300986055,1.9203374153835,1.919926519294563,1.919704622555685,1.918889217428244,1.91879563971699,1.918891088489886,1.918949027330683,1.918790148281943,1.920041000511134,1.922191335121028,1.920921534994307,1.918792166101838,1.922199966746405)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2+$i, 15).Value = $values[$i]
}

